$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A4: new date-formatted cell (numFmtId 16 "d-mmm"), new fill (white) ---
$ws.Range("A4").Value = 39904.166666666664
$ws.Range("A4").Interior.Color = 16777215
$ws.Range("A4").NumberFormat = "d-mmm"

# --- A5: new date-formatted cell (numFmtId 14 "mm-dd-yy"), new fill (orange) ---
$ws.Range("A5").Value = 39904.166666666664
$ws.Range("A5").Interior.Color = 43775
$ws.Range("A5").NumberFormat = "mm-dd-yy"
